$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row of tracked data (21 Aug 2019, 12:10, 74 pages, 22658 words),
# continuing the log below the existing last row (50).
$ws.Range("A50").Copy()
$ws.Range("A51").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B50").Copy()
$ws.Range("B51").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A51").Value = 43698
$ws.Range("B51").Value = 0.50694444444444442
$ws.Range("C51").Value = 74
$ws.Range("D51").Value = 22658

$excel.CutCopyMode = $false

$ws.Range("B52").Select()
